$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J, copying H1's formatting then overwriting text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Fill data rows 2-31: I column is always 1, J column mirrors H column value
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
